# Update the "trading data" worksheet with the new UNG put trade values,
# replacing the previous TZA call trade values. Column A (field names) is
# unchanged; only column B (values) is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B values, row 2 through row 30 (row 1 is the FIELD/VALUE header)
$ws.Range("B2").Value  = "20160928 +UNG-161021P8.00"
$ws.Range("B3").Value  = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20160928\\20160928"
$ws.Range("B4").Value  = "put"
$ws.Range("B5").Value  = 8
$ws.Range("B6").Value  = 2016
$ws.Range("B7").Value  = 10
$ws.Range("B8").Value  = 21
$ws.Range("B9").Value  = 0.13
$ws.Range("B10").Value = 8.515
$ws.Range("B11").Value = 2016
$ws.Range("B12").Value = 9
$ws.Range("B13").Value = 28
$ws.Range("B14").Value = 10
$ws.Range("B15").Value = 46
$ws.Range("B16").Value = 14
$ws.Range("B17").Value = "UNG"
$ws.Range("B18").Value = 0.3085
$ws.Range("B19").Value = "20160928 +UNG-161021P9.00"
$ws.Range("B20").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20160928\\20160928"
$ws.Range("B21").Value = 9
$ws.Range("B22").Value = -0.65
$ws.Range("B23").Value = 8.465
$ws.Range("B24").Value = 2016
$ws.Range("B25").Value = 9
$ws.Range("B26").Value = 28
$ws.Range("B27").Value = 13
$ws.Range("B28").Value = 24
$ws.Range("B29").Value = 38
$ws.Range("B30").Value = 0.3085

# Move the active selection to B31 (below the data) to match the saved view state.
$ws.Range("B31").Select()
